$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Accuracy" (sheet1): add a new "TDD Generated Tests" column (E)
# ---------------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("Accuracy")

$wsAcc.Range("E1").Value = "TDD Generated Tests"

$wsAcc.Range("E2").Value = 39.6
$wsAcc.Range("E2").NumberFormat = "0.00"

$wsAcc.Range("E3").Value = 40.6

$wsAcc.Range("E4").Value = 44.8
$wsAcc.Range("E4").NumberFormat = "0.00"

$wsAcc.Range("E5").Value = 51.6
$wsAcc.Range("E5").NumberFormat = "0.00"

$wsAcc.Range("E6").Value = 54.6
$wsAcc.Range("E6").NumberFormat = "0.00"

$wsAcc.Range("E7").Value = 50.8
$wsAcc.Range("E7").NumberFormat = "0.00"

$wsAcc.Range("E8").Value = 52.2
$wsAcc.Range("E8").NumberFormat = "0.00"

# widen the new column to match the existing "E" column width
$wsAcc.Columns.Item(6).ColumnWidth = 17.8

# ---------------------------------------------------------------------------
# Sheet "Errors" (sheet2): add a new merged "TDD Generated Tests" block
# (columns K:L:M) mirroring the existing Vanilla API / Google Translate /
# NLLB-200 Translate blocks.
# ---------------------------------------------------------------------------
$wsErr = $wb.Worksheets.Item("Errors")

$wsErr.Range("K1:M1").HorizontalAlignment = -4108
$wsErr.Range("K1").Value = "TDD Generated Tests"
$wsErr.Range("K1:M1").Merge()

$wsErr.Range("K2").Value = "Assertion"
$wsErr.Range("L2").Value = "Runtime"
$wsErr.Range("M2").Value = "Compilation"

$wsErr.Range("K3").Value = 214
$wsErr.Range("L3").Value = 46
$wsErr.Range("M3").Value = 42

$wsErr.Range("K4").Value = 230
$wsErr.Range("L4").Value = 41
$wsErr.Range("M4").Value = 26

$wsErr.Range("K5").Value = 234
$wsErr.Range("L5").Value = 41
$wsErr.Range("M5").Value = 1

$wsErr.Range("K6").Value = 196
$wsErr.Range("L6").Value = 28
$wsErr.Range("M6").Value = 18

$wsErr.Range("K7").Value = 191
$wsErr.Range("L7").Value = 30
$wsErr.Range("M7").Value = 6

$wsErr.Range("K8").Value = 170
$wsErr.Range("L8").Value = 67
$wsErr.Range("M8").Value = 9

$wsErr.Range("K9").Value = 188
$wsErr.Range("L9").Value = 49
$wsErr.Range("M9").Value = 2

# ---------------------------------------------------------------------------
# Selections / active sheet: the "Errors" tab becomes the active tab, with
# a new selection on each sheet.
# ---------------------------------------------------------------------------
$wsAcc.Range("F13").Select()
$wsErr.Activate()
$wsErr.Range("L19").Select()
